$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update "last updated" timestamp text (A1) ---
$ws.Range("A1").Value = "Datos actualizados a 21 de Abril de 2020 a las 06:52"

# --- Row 34 (Mexico): only "Nuevos casos" (F) column updated ---
$ws.Range("F34").Value = 378

# --- Countries reordered: Australia now listed before Serbia ---
# Row 41 becomes Australia (fresh data), Row 42 becomes Serbia (previous row 41 data)
$ws.Range("A41").Value = "Australia"
$ws.Range("B41").Value = 6642
$ws.Range("C41").Value = 17
$ws.Range("D41").Value = 4685
$ws.Range("E41").Value = 1886
$ws.Range("F41").Value = 49
$ws.Range("G41").Value = 0
$ws.Range("H41").Value = 71

$ws.Range("A42").Value = "Serbia"
$ws.Range("B42").Value = 6630
$ws.Range("C42").Value = 0
$ws.Range("D42").Value = 870
$ws.Range("E42").Value = 5635
$ws.Range("F42").Value = 108
$ws.Range("G42").Value = 0
$ws.Range("H42").Value = 125

# --- Row 58 (Tailandia): data refresh ---
$ws.Range("B58").Value = 2811
$ws.Range("C58").Value = 19
$ws.Range("D58").Value = 2108
$ws.Range("E58").Value = 655
$ws.Range("F58").Value = 61
$ws.Range("G58").Value = 1
$ws.Range("H58").Value = 48

# --- Countries reordered: Kirguistan now listed before Albania & Burkina Faso ---
# Row 99 becomes Kirguistan (fresh data)
# Row 100 becomes Albania (previous row 99 data)
# Row 101 becomes Burkina Faso (previous row 100 data)
$ws.Range("A99").Value = "Kirguistan"
$ws.Range("B99").Value = 590
$ws.Range("C99").Value = 22
$ws.Range("D99").Value = 216
$ws.Range("E99").Value = 367
$ws.Range("F99").Value = 5
$ws.Range("G99").Value = 0
$ws.Range("H99").Value = 7

$ws.Range("A100").Value = "Albania"
$ws.Range("B100").Value = 584
$ws.Range("C100").Value = 0
$ws.Range("D100").Value = 327
$ws.Range("E100").Value = 231
$ws.Range("F100").Value = 5
$ws.Range("G100").Value = 0
$ws.Range("H100").Value = 26

$ws.Range("A101").Value = "Burkina Faso"
$ws.Range("B101").Value = 581
$ws.Range("C101").Value = 0
$ws.Range("D101").Value = 357
$ws.Range("E101").Value = 186
$ws.Range("F101").Value = 0
$ws.Range("G101").Value = 0
$ws.Range("H101").Value = 38

# --- Row 164 (Macao): data refresh ---
$ws.Range("D164").Value = 24
$ws.Range("E164").Value = 21

# --- Countries reordered: Mongolia now listed before Republica del Chad ---
# Row 171 becomes Mongolia (fresh data)
# Row 172 becomes Republica del Chad (previous row 171 data)
$ws.Range("A171").Value = "Mongolia"
$ws.Range("B171").Value = 34
$ws.Range("C171").Value = 1
$ws.Range("D171").Value = 8
$ws.Range("E171").Value = 26
$ws.Range("F171").Value = 0
$ws.Range("G171").Value = 0
$ws.Range("H171").Value = 0

$ws.Range("A172").Value = "Republica del Chad"
$ws.Range("B172").Value = 33
$ws.Range("C172").Value = 0
$ws.Range("D172").Value = 8
$ws.Range("E172").Value = 25
$ws.Range("F172").Value = 0
$ws.Range("G172").Value = 0
$ws.Range("H172").Value = 0
